$wb = $excel.ActiveWorkbook

$sheetNames = @("LP1912", "6203-6173")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Insert a new row at position 12, shifting the existing row 12 (and below) down to 13
    $ws.Rows.Item(12).Insert()

    # Fill in the new row 12 with the new schedule entry
    $ws.Cells.Item(12, 1).Value = "03:00:53"
    $ws.Cells.Item(12, 2).Value = "04:53"
    $ws.Cells.Item(12, 3).Value = "11_ETCHEVERRY"
    $ws.Cells.Item(12, 4).Value = 113

    # Update header info
    $ws.Range("A2").Value = "Última actualización: 03:00:53"
    $ws.Range("A3").Value = "Total filas: 8"
}
